$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph - this is the anchor for the
# block of boilerplate (footer) text that was stripped from the page in
# this revision.
$jupiter = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $jupiter = $p
        break
    }
}

if ($jupiter -ne $null) {
    # The paragraph immediately before it is the blank spacer paragraph
    # that separated it from "LOB1019: Física II (Requisito)", and the
    # paragraph immediately after it is the "© 2020 ..." copyright line.
    # Both the blank spacer and the copyright paragraph are removed along
    # with the "Ver no Jupiter ..." paragraph itself, leaving the single
    # blank paragraph that precedes the trailing page break intact.
    $before = $jupiter.Previous()
    $after = $jupiter.Next()

    $start = $before.Range.Start
    $end = $after.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
